$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for B2:X5 (rows = data rows 2-5, cols = B..X)
$data = @(
    @(0.0053680981595092,0.00230061349693252,0.0498466257668712,0.031441717791411,0.00230061349693252,0.00920245398773006,0.00920245398773006,0.00230061349693252,0.00920245398773006,0.996165644171779,0.000766871165644172,0,0.000766871165644172,0.00306748466257669,0.00153374233128834,0.00460122699386503,0.941717791411043,0.0107361963190184,0.00460122699386503,0.00996932515337423,0.849693251533742,0.159509202453988,0.0368098159509202),
    @(0.976993865030675,0.973159509202454,0.00613496932515337,0.953220858895706,0.0184049079754601,0,0,0,0,0,0.980061349693252,0.00153374233128834,0.00460122699386503,0.00383435582822086,0.997699386503067,0.977760736196319,0,0.000766871165644172,0.0191717791411043,0.0176380368098159,0.0145705521472393,0.000766871165644172,0.0237730061349693),
    @(0.0153374233128834,0.0191717791411043,0.943251533742331,0,0.00230061349693252,0.989263803680982,0.989263803680982,0.996932515337423,0.990030674846626,0.00383435582822086,0.0153374233128834,0,0.0145705521472393,0.0153374233128834,0,0.0138036809815951,0.0575153374233129,0.983128834355828,0.975460122699387,0.968558282208589,0.131901840490798,0.825153374233129,0.93941717791411),
    @(0.00230061349693252,0.0053680981595092,0.000766871165644172,0.0153374233128834,0.976993865030675,0.000766871165644172,0.000766871165644172,0,0.000766871165644172,0,0.00383435582822086,0.998466257668712,0.980061349693252,0.977760736196319,0.000766871165644172,0.00383435582822086,0.000766871165644172,0.0053680981595092,0.000766871165644172,0.00383435582822086,0.00383435582822086,0.0145705521472393,0)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}

Write-Output "Updated B2:X5 with new frequency values"